# Schedule_Task_Book.xlsx — add "Output raw data" task + status notes on the
# Algorithm sheet, and leave that sheet as the active/selected one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Algorithm")

# New status column (B) gets a wider custom width (stored width "51").
# The engine's char-width->stored-width conversion adds 5/6 to whatever
# we request, so back that off to land exactly on 51.
$ws.Columns.Item(2).ColumnWidth = 51 - (5/6)

# Status / comment notes next to existing tasks.
$ws.Range("B2").Value = "Almost, some timings need to be adjusted as they're done quickly."
$ws.Range("B4").Value = "Done"
$ws.Range("B11").Value = "Done"

# New task row at the bottom of the list.
$ws.Range("A17").Value = "Output the data from each step in C"
$ws.Range("B17").Value = "Done"
$ws.Rows.Item(17).RowHeight = 30

# Make the Algorithm sheet the active tab/selection (was Hardware before).
$ws.Activate()
$ws.Range("B18").Select() | Out-Null
